$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.177.25"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.277.25"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "111.68"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "264.68"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.645"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.33%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.610"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.58"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0936"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.31"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.22%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.27"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.617.76"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.277.75"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.202.57"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.74"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.16"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.42"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.19"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.88"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.35"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.01%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.35"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "40.93"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.10%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.26%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.17"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.44"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0897"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.17%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.38%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0381"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.67%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.88"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.44%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.58"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +6.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "14.21"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "74.94"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.236"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.09%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.55%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.55"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.26%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.26"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0992"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "99.97"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.599"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +9.90%  "
